$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full re-shuffle of the data rows (2-131) plus one new weekly record appended as row 132.
$data = New-Object 'object[,]' 131,18
$data[0,0] = 9
$data[0,1] = 'Vega Central Mapocho de Santiago'
$data[0,2] = 'Metropolitana'
$data[0,3] = 44278
$data[0,4] = 13
$data[0,5] = 100112005
$data[0,6] = 'Puerro'
$data[0,7] = 'Sin especificar'
$data[0,8] = 'Primera'
$data[0,9] = 70
$data[0,10] = 8000
$data[0,11] = 8000
$data[0,12] = 8000
$data[0,13] = '$/paquete 20 unidades'
$data[0,14] = 'Provincia de Chacabuco'
$data[0,15] = 400
$data[0,16] = 20
$data[0,17] = 'Hortaliza'
$data[1,0] = 9
$data[1,1] = 'Vega Central Mapocho de Santiago'
$data[1,2] = 'Metropolitana'
$data[1,3] = 44468
$data[1,4] = 13
$data[1,5] = 100112005
$data[1,6] = 'Puerro'
$data[1,7] = 'Sin especificar'
$data[1,8] = 'Primera'
$data[1,9] = 133
$data[1,10] = 7000
$data[1,11] = 8000
$data[1,12] = 7504
$data[1,13] = '$/paquete 20 unidades'
$data[1,14] = 'Provincia de Chacabuco'
$data[1,15] = 375
$data[1,16] = 20
$data[1,17] = 'Hortaliza'
$data[2,0] = 9
$data[2,1] = 'Vega Central Mapocho de Santiago'
$data[2,2] = 'Metropolitana'
$data[2,3] = 44239
$data[2,4] = 13
$data[2,5] = 100112005
$data[2,6] = 'Puerro'
$data[2,7] = 'Sin especificar'
$data[2,8] = 'Primera'
$data[2,9] = 70
$data[2,10] = 8000
$data[2,11] = 8000
$data[2,12] = 8000
$data[2,13] = '$/paquete 20 unidades'
$data[2,14] = 'Provincia de Chacabuco'
$data[2,15] = 400
$data[2,16] = 20
$data[2,17] = 'Hortaliza'
$data[3,0] = 9
$data[3,1] = 'Vega Central Mapocho de Santiago'
$data[3,2] = 'Metropolitana'
$data[3,3] = 44539
$data[3,4] = 13
$data[3,5] = 100112005
$data[3,6] = 'Puerro'
$data[3,7] = 'Sin especificar'
$data[3,8] = 'Primera'
$data[3,9] = 133
$data[3,10] = 6000
$data[3,11] = 7000
$data[3,12] = 6504
$data[3,13] = '$/paquete 20 unidades'
$data[3,14] = 'Provincia de Chacabuco'
$data[3,15] = 325
$data[3,16] = 20
$data[3,17] = 'Hortaliza'
$data[4,0] = 9
$data[4,1] = 'Vega Central Mapocho de Santiago'
$data[4,2] = 'Metropolitana'
$data[4,3] = 44356
$data[4,4] = 13
$data[4,5] = 100112005
$data[4,6] = 'Puerro'
$data[4,7] = 'Sin especificar'
$data[4,8] = 'Primera'
$data[4,9] = 160
$data[4,10] = 7000
$data[4,11] = 8000
$data[4,12] = 7500
$data[4,13] = '$/paquete 20 unidades'
$data[4,14] = 'Provincia de Chacabuco'
$data[4,15] = 375
$data[4,16] = 20
$data[4,17] = 'Hortaliza'
$data[5,0] = 9
$data[5,1] = 'Vega Central Mapocho de Santiago'
$data[5,2] = 'Metropolitana'
$data[5,3] = 44671
$data[5,4] = 13
$data[5,5] = 100112005
$data[5,6] = 'Puerro'
$data[5,7] = 'Sin especificar'
$data[5,8] = 'Primera'
$data[5,9] = 160
$data[5,10] = 8000
$data[5,11] = 8000
$data[5,12] = 8000
$data[5,13] = '$/paquete 20 unidades'
$data[5,14] = 'Provincia de Chacabuco'
$data[5,15] = 400
$data[5,16] = 20
$data[5,17] = 'Hortaliza'
$data[6,0] = 9
$data[6,1] = 'Vega Central Mapocho de Santiago'
$data[6,2] = 'Metropolitana'
$data[6,3] = 44414
$data[6,4] = 13
$data[6,5] = 100112005
$data[6,6] = 'Puerro'
$data[6,7] = 'Sin especificar'
$data[6,8] = 'Primera'
$data[6,9] = 180
$data[6,10] = 7500
$data[6,11] = 8000
$data[6,12] = 7750
$data[6,13] = '$/paquete 20 unidades'
$data[6,14] = 'Provincia de Chacabuco'
$data[6,15] = 388
$data[6,16] = 20
$data[6,17] = 'Hortaliza'
$data[7,0] = 9
$data[7,1] = 'Vega Central Mapocho de Santiago'
$data[7,2] = 'Metropolitana'
$data[7,3] = 44384
$data[7,4] = 13
$data[7,5] = 100112005
$data[7,6] = 'Puerro'
$data[7,7] = 'Sin especificar'
$data[7,8] = 'Primera'
$data[7,9] = 160
$data[7,10] = 8000
$data[7,11] = 9000
$data[7,12] = 8500
$data[7,13] = '$/paquete 20 unidades'
$data[7,14] = 'Provincia de Chacabuco'
$data[7,15] = 425
$data[7,16] = 20
$data[7,17] = 'Hortaliza'
$data[8,0] = 9
$data[8,1] = 'Vega Central Mapocho de Santiago'
$data[8,2] = 'Metropolitana'
$data[8,3] = 44308
$data[8,4] = 13
$data[8,5] = 100112005
$data[8,6] = 'Puerro'
$data[8,7] = 'Sin especificar'
$data[8,8] = 'Primera'
$data[8,9] = 160
$data[8,10] = 8000
$data[8,11] = 8000
$data[8,12] = 8000
$data[8,13] = '$/paquete 20 unidades'
$data[8,14] = 'Provincia de Chacabuco'
$data[8,15] = 400
$data[8,16] = 20
$data[8,17] = 'Hortaliza'
$data[9,0] = 9
$data[9,1] = 'Vega Central Mapocho de Santiago'
$data[9,2] = 'Metropolitana'
$data[9,3] = 44517
$data[9,4] = 13
$data[9,5] = 100112005
$data[9,6] = 'Puerro'
$data[9,7] = 'Sin especificar'
$data[9,8] = 'Primera'
$data[9,9] = 160
$data[9,10] = 6000
$data[9,11] = 7000
$data[9,12] = 6500
$data[9,13] = '$/paquete 20 unidades'
$data[9,14] = 'Provincia de Chacabuco'
$data[9,15] = 325
$data[9,16] = 20
$data[9,17] = 'Hortaliza'
$data[10,0] = 9
$data[10,1] = 'Vega Central Mapocho de Santiago'
$data[10,2] = 'Metropolitana'
$data[10,3] = 44344
$data[10,4] = 13
$data[10,5] = 100112005
$data[10,6] = 'Puerro'
$data[10,7] = 'Sin especificar'
$data[10,8] = 'Primera'
$data[10,9] = 210
$data[10,10] = 8000
$data[10,11] = 8000
$data[10,12] = 8000
$data[10,13] = '$/paquete 20 unidades'
$data[10,14] = 'Provincia de Chacabuco'
$data[10,15] = 400
$data[10,16] = 20
$data[10,17] = 'Hortaliza'
$data[11,0] = 9
$data[11,1] = 'Vega Central Mapocho de Santiago'
$data[11,2] = 'Metropolitana'
$data[11,3] = 44391
$data[11,4] = 13
$data[11,5] = 100112005
$data[11,6] = 'Puerro'
$data[11,7] = 'Sin especificar'
$data[11,8] = 'Primera'
$data[11,9] = 52
$data[11,10] = 7000
$data[11,11] = 8000
$data[11,12] = 7500
$data[11,13] = '$/paquete 20 unidades'
$data[11,14] = 'Provincia de Chacabuco'
$data[11,15] = 375
$data[11,16] = 20
$data[11,17] = 'Hortaliza'
$data[12,0] = 9
$data[12,1] = 'Vega Central Mapocho de Santiago'
$data[12,2] = 'Metropolitana'
$data[12,3] = 44475
$data[12,4] = 13
$data[12,5] = 100112005
$data[12,6] = 'Puerro'
$data[12,7] = 'Sin especificar'
$data[12,8] = 'Primera'
$data[12,9] = 160
$data[12,10] = 7000
$data[12,11] = 8000
$data[12,12] = 7500
$data[12,13] = '$/paquete 20 unidades'
$data[12,14] = 'Provincia de Chacabuco'
$data[12,15] = 375
$data[12,16] = 20
$data[12,17] = 'Hortaliza'
$data[13,0] = 9
$data[13,1] = 'Vega Central Mapocho de Santiago'
$data[13,2] = 'Metropolitana'
$data[13,3] = 45112
$data[13,4] = 13
$data[13,5] = 100112005
$data[13,6] = 'Puerro'
$data[13,7] = 'Sin especificar'
$data[13,8] = 'Primera'
$data[13,9] = 70
$data[13,10] = 8000
$data[13,11] = 8000
$data[13,12] = 8000
$data[13,13] = '$/paquete 20 unidades'
$data[13,14] = 'Provincia de Chacabuco'
$data[13,15] = 400
$data[13,16] = 20
$data[13,17] = 'Hortaliza'
$data[14,0] = 9
$data[14,1] = 'Vega Central Mapocho de Santiago'
$data[14,2] = 'Metropolitana'
$data[14,3] = 44174
$data[14,4] = 13
$data[14,5] = 100112005
$data[14,6] = 'Puerro'
$data[14,7] = 'Sin especificar'
$data[14,8] = 'Primera'
$data[14,9] = 70
$data[14,10] = 8000
$data[14,11] = 8000
$data[14,12] = 8000
$data[14,13] = '$/paquete 20 unidades'
$data[14,14] = 'Provincia de Chacabuco'
$data[14,15] = 400
$data[14,16] = 20
$data[14,17] = 'Hortaliza'
$data[15,0] = 9
$data[15,1] = 'Vega Central Mapocho de Santiago'
$data[15,2] = 'Metropolitana'
$data[15,3] = 44769
$data[15,4] = 13
$data[15,5] = 100112005
$data[15,6] = 'Puerro'
$data[15,7] = 'Sin especificar'
$data[15,8] = 'Primera'
$data[15,9] = 70
$data[15,10] = 7000
$data[15,11] = 8000
$data[15,12] = 7500
$data[15,13] = '$/paquete 20 unidades'
$data[15,14] = 'Provincia de Chacabuco'
$data[15,15] = 375
$data[15,16] = 20
$data[15,17] = 'Hortaliza'
$data[16,0] = 9
$data[16,1] = 'Vega Central Mapocho de Santiago'
$data[16,2] = 'Metropolitana'
$data[16,3] = 44265
$data[16,4] = 13
$data[16,5] = 100112005
$data[16,6] = 'Puerro'
$data[16,7] = 'Sin especificar'
$data[16,8] = 'Primera'
$data[16,9] = 70
$data[16,10] = 8000
$data[16,11] = 8000
$data[16,12] = 8000
$data[16,13] = '$/paquete 20 unidades'
$data[16,14] = 'Provincia de Chacabuco'
$data[16,15] = 400
$data[16,16] = 20
$data[16,17] = 'Hortaliza'
$data[17,0] = 9
$data[17,1] = 'Vega Central Mapocho de Santiago'
$data[17,2] = 'Metropolitana'
$data[17,3] = 44519
$data[17,4] = 13
$data[17,5] = 100112005
$data[17,6] = 'Puerro'
$data[17,7] = 'Sin especificar'
$data[17,8] = 'Primera'
$data[17,9] = 160
$data[17,10] = 6000
$data[17,11] = 7000
$data[17,12] = 6500
$data[17,13] = '$/paquete 20 unidades'
$data[17,14] = 'Provincia de Chacabuco'
$data[17,15] = 325
$data[17,16] = 20
$data[17,17] = 'Hortaliza'
$data[18,0] = 9
$data[18,1] = 'Vega Central Mapocho de Santiago'
$data[18,2] = 'Metropolitana'
$data[18,3] = 44208
$data[18,4] = 13
$data[18,5] = 100112005
$data[18,6] = 'Puerro'
$data[18,7] = 'Sin especificar'
$data[18,8] = 'Primera'
$data[18,9] = 50
$data[18,10] = 8000
$data[18,11] = 8000
$data[18,12] = 8000
$data[18,13] = '$/paquete 20 unidades'
$data[18,14] = 'Provincia de Chacabuco'
$data[18,15] = 400
$data[18,16] = 20
$data[18,17] = 'Hortaliza'
$data[19,0] = 9
$data[19,1] = 'Vega Central Mapocho de Santiago'
$data[19,2] = 'Metropolitana'
$data[19,3] = 44167
$data[19,4] = 13
$data[19,5] = 100112005
$data[19,6] = 'Puerro'
$data[19,7] = 'Sin especificar'
$data[19,8] = 'Primera'
$data[19,9] = 50
$data[19,10] = 8000
$data[19,11] = 8000
$data[19,12] = 8000
$data[19,13] = '$/paquete 20 unidades'
$data[19,14] = 'Provincia de Chacabuco'
$data[19,15] = 400
$data[19,16] = 20
$data[19,17] = 'Hortaliza'
$data[20,0] = 9
$data[20,1] = 'Vega Central Mapocho de Santiago'
$data[20,2] = 'Metropolitana'
$data[20,3] = 44755
$data[20,4] = 13
$data[20,5] = 100112005
$data[20,6] = 'Puerro'
$data[20,7] = 'Sin especificar'
$data[20,8] = 'Primera'
$data[20,9] = 160
$data[20,10] = 7000
$data[20,11] = 8000
$data[20,12] = 7500
$data[20,13] = '$/paquete 20 unidades'
$data[20,14] = 'Provincia de Chacabuco'
$data[20,15] = 375
$data[20,16] = 20
$data[20,17] = 'Hortaliza'
$data[21,0] = 9
$data[21,1] = 'Vega Central Mapocho de Santiago'
$data[21,2] = 'Metropolitana'
$data[21,3] = 44421
$data[21,4] = 13
$data[21,5] = 100112005
$data[21,6] = 'Puerro'
$data[21,7] = 'Sin especificar'
$data[21,8] = 'Primera'
$data[21,9] = 180
$data[21,10] = 7000
$data[21,11] = 8000
$data[21,12] = 7500
$data[21,13] = '$/paquete 20 unidades'
$data[21,14] = 'Provincia de Chacabuco'
$data[21,15] = 375
$data[21,16] = 20
$data[21,17] = 'Hortaliza'
$data[22,0] = 9
$data[22,1] = 'Vega Central Mapocho de Santiago'
$data[22,2] = 'Metropolitana'
$data[22,3] = 44349
$data[22,4] = 13
$data[22,5] = 100112005
$data[22,6] = 'Puerro'
$data[22,7] = 'Sin especificar'
$data[22,8] = 'Primera'
$data[22,9] = 130
$data[22,10] = 8000
$data[22,11] = 8000
$data[22,12] = 8000
$data[22,13] = '$/paquete 20 unidades'
$data[22,14] = 'Provincia de Chacabuco'
$data[22,15] = 400
$data[22,16] = 20
$data[22,17] = 'Hortaliza'
$data[23,0] = 9
$data[23,1] = 'Vega Central Mapocho de Santiago'
$data[23,2] = 'Metropolitana'
$data[23,3] = 44224
$data[23,4] = 13
$data[23,5] = 100112005
$data[23,6] = 'Puerro'
$data[23,7] = 'Sin especificar'
$data[23,8] = 'Primera'
$data[23,9] = 120
$data[23,10] = 6000
$data[23,11] = 7000
$data[23,12] = 6667
$data[23,13] = '$/paquete 20 unidades'
$data[23,14] = 'Provincia de Chacabuco'
$data[23,15] = 333
$data[23,16] = 20
$data[23,17] = 'Hortaliza'
$data[24,0] = 9
$data[24,1] = 'Vega Central Mapocho de Santiago'
$data[24,2] = 'Metropolitana'
$data[24,3] = 44365
$data[24,4] = 13
$data[24,5] = 100112005
$data[24,6] = 'Puerro'
$data[24,7] = 'Sin especificar'
$data[24,8] = 'Primera'
$data[24,9] = 180
$data[24,10] = 8000
$data[24,11] = 8000
$data[24,12] = 8000
$data[24,13] = '$/paquete 20 unidades'
$data[24,14] = 'Provincia de Chacabuco'
$data[24,15] = 400
$data[24,16] = 20
$data[24,17] = 'Hortaliza'
$data[25,0] = 9
$data[25,1] = 'Vega Central Mapocho de Santiago'
$data[25,2] = 'Metropolitana'
$data[25,3] = 44748
$data[25,4] = 13
$data[25,5] = 100112005
$data[25,6] = 'Puerro'
$data[25,7] = 'Sin especificar'
$data[25,8] = 'Primera'
$data[25,9] = 70
$data[25,10] = 8000
$data[25,11] = 8000
$data[25,12] = 8000
$data[25,13] = '$/paquete 20 unidades'
$data[25,14] = 'Provincia de Chacabuco'
$data[25,15] = 400
$data[25,16] = 20
$data[25,17] = 'Hortaliza'
$data[26,0] = 9
$data[26,1] = 'Vega Central Mapocho de Santiago'
$data[26,2] = 'Metropolitana'
$data[26,3] = 44657
$data[26,4] = 13
$data[26,5] = 100112005
$data[26,6] = 'Puerro'
$data[26,7] = 'Sin especificar'
$data[26,8] = 'Primera'
$data[26,9] = 160
$data[26,10] = 7000
$data[26,11] = 8000
$data[26,12] = 7500
$data[26,13] = '$/paquete 20 unidades'
$data[26,14] = 'Provincia de Chacabuco'
$data[26,15] = 375
$data[26,16] = 20
$data[26,17] = 'Hortaliza'
$data[27,0] = 9
$data[27,1] = 'Vega Central Mapocho de Santiago'
$data[27,2] = 'Metropolitana'
$data[27,3] = 44643
$data[27,4] = 13
$data[27,5] = 100112005
$data[27,6] = 'Puerro'
$data[27,7] = 'Sin especificar'
$data[27,8] = 'Primera'
$data[27,9] = 160
$data[27,10] = 8000
$data[27,11] = 9000
$data[27,12] = 8500
$data[27,13] = '$/paquete 20 unidades'
$data[27,14] = 'Provincia de Chacabuco'
$data[27,15] = 425
$data[27,16] = 20
$data[27,17] = 'Hortaliza'
$data[28,0] = 9
$data[28,1] = 'Vega Central Mapocho de Santiago'
$data[28,2] = 'Metropolitana'
$data[28,3] = 44477
$data[28,4] = 13
$data[28,5] = 100112005
$data[28,6] = 'Puerro'
$data[28,7] = 'Sin especificar'
$data[28,8] = 'Primera'
$data[28,9] = 160
$data[28,10] = 7000
$data[28,11] = 8000
$data[28,12] = 7500
$data[28,13] = '$/paquete 20 unidades'
$data[28,14] = 'Provincia de Chacabuco'
$data[28,15] = 375
$data[28,16] = 20
$data[28,17] = 'Hortaliza'
$data[29,0] = 9
$data[29,1] = 'Vega Central Mapocho de Santiago'
$data[29,2] = 'Metropolitana'
$data[29,3] = 45070
$data[29,4] = 13
$data[29,5] = 100112005
$data[29,6] = 'Puerro'
$data[29,7] = 'Sin especificar'
$data[29,8] = 'Primera'
$data[29,9] = 70
$data[29,10] = 7000
$data[29,11] = 7000
$data[29,12] = 7000
$data[29,13] = '$/paquete 20 unidades'
$data[29,14] = 'Provincia de Chacabuco'
$data[29,15] = 350
$data[29,16] = 20
$data[29,17] = 'Hortaliza'
$data[30,0] = 9
$data[30,1] = 'Vega Central Mapocho de Santiago'
$data[30,2] = 'Metropolitana'
$data[30,3] = 44428
$data[30,4] = 13
$data[30,5] = 100112005
$data[30,6] = 'Puerro'
$data[30,7] = 'Sin especificar'
$data[30,8] = 'Primera'
$data[30,9] = 97
$data[30,10] = 8000
$data[30,11] = 9000
$data[30,12] = 8505
$data[30,13] = '$/paquete 20 unidades'
$data[30,14] = 'Provincia de Chacabuco'
$data[30,15] = 425
$data[30,16] = 20
$data[30,17] = 'Hortaliza'
$data[31,0] = 9
$data[31,1] = 'Vega Central Mapocho de Santiago'
$data[31,2] = 'Metropolitana'
$data[31,3] = 45126
$data[31,4] = 13
$data[31,5] = 100112005
$data[31,6] = 'Puerro'
$data[31,7] = 'Sin especificar'
$data[31,8] = 'Primera'
$data[31,9] = 70
$data[31,10] = 8000
$data[31,11] = 8000
$data[31,12] = 8000
$data[31,13] = '$/paquete 20 unidades'
$data[31,14] = 'Provincia de Chacabuco'
$data[31,15] = 400
$data[31,16] = 20
$data[31,17] = 'Hortaliza'
$data[32,0] = 9
$data[32,1] = 'Vega Central Mapocho de Santiago'
$data[32,2] = 'Metropolitana'
$data[32,3] = 44335
$data[32,4] = 13
$data[32,5] = 100112005
$data[32,6] = 'Puerro'
$data[32,7] = 'Sin especificar'
$data[32,8] = 'Primera'
$data[32,9] = 250
$data[32,10] = 7000
$data[32,11] = 8000
$data[32,12] = 7500
$data[32,13] = '$/paquete 20 unidades'
$data[32,14] = 'Provincia de Chacabuco'
$data[32,15] = 375
$data[32,16] = 20
$data[32,17] = 'Hortaliza'
$data[33,0] = 9
$data[33,1] = 'Vega Central Mapocho de Santiago'
$data[33,2] = 'Metropolitana'
$data[33,3] = 44498
$data[33,4] = 13
$data[33,5] = 100112005
$data[33,6] = 'Puerro'
$data[33,7] = 'Sin especificar'
$data[33,8] = 'Primera'
$data[33,9] = 130
$data[33,10] = 9000
$data[33,11] = 10000
$data[33,12] = 9615
$data[33,13] = '$/paquete 20 unidades'
$data[33,14] = 'Provincia de Santiago'
$data[33,15] = 481
$data[33,16] = 20
$data[33,17] = 'Hortaliza'
$data[34,0] = 9
$data[34,1] = 'Vega Central Mapocho de Santiago'
$data[34,2] = 'Metropolitana'
$data[34,3] = 44679
$data[34,4] = 13
$data[34,5] = 100112005
$data[34,6] = 'Puerro'
$data[34,7] = 'Sin especificar'
$data[34,8] = 'Primera'
$data[34,9] = 97
$data[34,10] = 8000
$data[34,11] = 9000
$data[34,12] = 8505
$data[34,13] = '$/paquete 20 unidades'
$data[34,14] = 'Provincia de Chacabuco'
$data[34,15] = 425
$data[34,16] = 20
$data[34,17] = 'Hortaliza'
$data[35,0] = 9
$data[35,1] = 'Vega Central Mapocho de Santiago'
$data[35,2] = 'Metropolitana'
$data[35,3] = 45049
$data[35,4] = 13
$data[35,5] = 100112005
$data[35,6] = 'Puerro'
$data[35,7] = 'Sin especificar'
$data[35,8] = 'Primera'
$data[35,9] = 160
$data[35,10] = 9000
$data[35,11] = 9000
$data[35,12] = 9000
$data[35,13] = '$/paquete 20 unidades'
$data[35,14] = 'Provincia de Chacabuco'
$data[35,15] = 450
$data[35,16] = 20
$data[35,17] = 'Hortaliza'
$data[36,0] = 9
$data[36,1] = 'Vega Central Mapocho de Santiago'
$data[36,2] = 'Metropolitana'
$data[36,3] = 44461
$data[36,4] = 13
$data[36,5] = 100112005
$data[36,6] = 'Puerro'
$data[36,7] = 'Sin especificar'
$data[36,8] = 'Primera'
$data[36,9] = 79
$data[36,10] = 7000
$data[36,11] = 8000
$data[36,12] = 7494
$data[36,13] = '$/paquete 20 unidades'
$data[36,14] = 'Provincia de Chacabuco'
$data[36,15] = 375
$data[36,16] = 20
$data[36,17] = 'Hortaliza'
$data[37,0] = 9
$data[37,1] = 'Vega Central Mapocho de Santiago'
$data[37,2] = 'Metropolitana'
$data[37,3] = 44463
$data[37,4] = 13
$data[37,5] = 100112005
$data[37,6] = 'Puerro'
$data[37,7] = 'Sin especificar'
$data[37,8] = 'Primera'
$data[37,9] = 160
$data[37,10] = 7500
$data[37,11] = 8000
$data[37,12] = 7750
$data[37,13] = '$/paquete 20 unidades'
$data[37,14] = 'Provincia de Chacabuco'
$data[37,15] = 388
$data[37,16] = 20
$data[37,17] = 'Hortaliza'
$data[38,0] = 9
$data[38,1] = 'Vega Central Mapocho de Santiago'
$data[38,2] = 'Metropolitana'
$data[38,3] = 44847
$data[38,4] = 13
$data[38,5] = 100112005
$data[38,6] = 'Puerro'
$data[38,7] = 'Sin especificar'
$data[38,8] = 'Primera'
$data[38,9] = 70
$data[38,10] = 11000
$data[38,11] = 12000
$data[38,12] = 11571
$data[38,13] = '$/paquete 20 unidades'
$data[38,14] = 'Provincia de Melipilla'
$data[38,15] = 579
$data[38,16] = 20
$data[38,17] = 'Hortaliza'
$data[39,0] = 9
$data[39,1] = 'Vega Central Mapocho de Santiago'
$data[39,2] = 'Metropolitana'
$data[39,3] = 44847
$data[39,4] = 13
$data[39,5] = 100112005
$data[39,6] = 'Puerro'
$data[39,7] = 'Sin especificar'
$data[39,8] = 'Segunda'
$data[39,9] = 40
$data[39,10] = 10000
$data[39,11] = 10000
$data[39,12] = 10000
$data[39,13] = '$/paquete 20 unidades'
$data[39,14] = 'Provincia de Melipilla'
$data[39,15] = 500
$data[39,16] = 20
$data[39,17] = 'Hortaliza'
$data[40,0] = 9
$data[40,1] = 'Vega Central Mapocho de Santiago'
$data[40,2] = 'Metropolitana'
$data[40,3] = 44630
$data[40,4] = 13
$data[40,5] = 100112005
$data[40,6] = 'Puerro'
$data[40,7] = 'Sin especificar'
$data[40,8] = 'Primera'
$data[40,9] = 79
$data[40,10] = 9000
$data[40,11] = 10000
$data[40,12] = 9494
$data[40,13] = '$/paquete 20 unidades'
$data[40,14] = 'Provincia de Chacabuco'
$data[40,15] = 475
$data[40,16] = 20
$data[40,17] = 'Hortaliza'
$data[41,0] = 9
$data[41,1] = 'Vega Central Mapocho de Santiago'
$data[41,2] = 'Metropolitana'
$data[41,3] = 44232
$data[41,4] = 13
$data[41,5] = 100112005
$data[41,6] = 'Puerro'
$data[41,7] = 'Sin especificar'
$data[41,8] = 'Primera'
$data[41,9] = 60
$data[41,10] = 7000
$data[41,11] = 7000
$data[41,12] = 7000
$data[41,13] = '$/paquete 20 unidades'
$data[41,14] = 'Provincia de Chacabuco'
$data[41,15] = 350
$data[41,16] = 20
$data[41,17] = 'Hortaliza'
$data[42,0] = 9
$data[42,1] = 'Vega Central Mapocho de Santiago'
$data[42,2] = 'Metropolitana'
$data[42,3] = 44286
$data[42,4] = 13
$data[42,5] = 100112005
$data[42,6] = 'Puerro'
$data[42,7] = 'Sin especificar'
$data[42,8] = 'Primera'
$data[42,9] = 70
$data[42,10] = 8000
$data[42,11] = 8000
$data[42,12] = 8000
$data[42,13] = '$/paquete 20 unidades'
$data[42,14] = 'Provincia de Chacabuco'
$data[42,15] = 400
$data[42,16] = 20
$data[42,17] = 'Hortaliza'
$data[43,0] = 9
$data[43,1] = 'Vega Central Mapocho de Santiago'
$data[43,2] = 'Metropolitana'
$data[43,3] = 44299
$data[43,4] = 13
$data[43,5] = 100112005
$data[43,6] = 'Puerro'
$data[43,7] = 'Sin especificar'
$data[43,8] = 'Primera'
$data[43,9] = 160
$data[43,10] = 8000
$data[43,11] = 8000
$data[43,12] = 8000
$data[43,13] = '$/paquete 20 unidades'
$data[43,14] = 'Provincia de Chacabuco'
$data[43,15] = 400
$data[43,16] = 20
$data[43,17] = 'Hortaliza'
$data[44,0] = 9
$data[44,1] = 'Vega Central Mapocho de Santiago'
$data[44,2] = 'Metropolitana'
$data[44,3] = 44301
$data[44,4] = 13
$data[44,5] = 100112005
$data[44,6] = 'Puerro'
$data[44,7] = 'Sin especificar'
$data[44,8] = 'Primera'
$data[44,9] = 160
$data[44,10] = 8000
$data[44,11] = 8000
$data[44,12] = 8000
$data[44,13] = '$/paquete 20 unidades'
$data[44,14] = 'Provincia de Chacabuco'
$data[44,15] = 400
$data[44,16] = 20
$data[44,17] = 'Hortaliza'
$data[45,0] = 9
$data[45,1] = 'Vega Central Mapocho de Santiago'
$data[45,2] = 'Metropolitana'
$data[45,3] = 44650
$data[45,4] = 13
$data[45,5] = 100112005
$data[45,6] = 'Puerro'
$data[45,7] = 'Sin especificar'
$data[45,8] = 'Primera'
$data[45,9] = 160
$data[45,10] = 9000
$data[45,11] = 10000
$data[45,12] = 9500
$data[45,13] = '$/paquete 20 unidades'
$data[45,14] = 'Provincia de Chacabuco'
$data[45,15] = 475
$data[45,16] = 20
$data[45,17] = 'Hortaliza'
$data[46,0] = 9
$data[46,1] = 'Vega Central Mapocho de Santiago'
$data[46,2] = 'Metropolitana'
$data[46,3] = 44708
$data[46,4] = 13
$data[46,5] = 100112005
$data[46,6] = 'Puerro'
$data[46,7] = 'Sin especificar'
$data[46,8] = 'Primera'
$data[46,9] = 80
$data[46,10] = 7000
$data[46,11] = 7000
$data[46,12] = 7000
$data[46,13] = '$/paquete 20 unidades'
$data[46,14] = 'Provincia de Chacabuco'
$data[46,15] = 350
$data[46,16] = 20
$data[46,17] = 'Hortaliza'
$data[47,0] = 9
$data[47,1] = 'Vega Central Mapocho de Santiago'
$data[47,2] = 'Metropolitana'
$data[47,3] = 44873
$data[47,4] = 13
$data[47,5] = 100112005
$data[47,6] = 'Puerro'
$data[47,7] = 'Sin especificar'
$data[47,8] = 'Primera'
$data[47,9] = 70
$data[47,10] = 8000
$data[47,11] = 9000
$data[47,12] = 8571
$data[47,13] = '$/paquete 20 unidades'
$data[47,14] = 'Provincia de Chacabuco'
$data[47,15] = 429
$data[47,16] = 20
$data[47,17] = 'Hortaliza'
$data[48,0] = 9
$data[48,1] = 'Vega Central Mapocho de Santiago'
$data[48,2] = 'Metropolitana'
$data[48,3] = 44195
$data[48,4] = 13
$data[48,5] = 100112005
$data[48,6] = 'Puerro'
$data[48,7] = 'Sin especificar'
$data[48,8] = 'Primera'
$data[48,9] = 70
$data[48,10] = 7000
$data[48,11] = 7000
$data[48,12] = 7000
$data[48,13] = '$/paquete 20 unidades'
$data[48,14] = 'Provincia de Chacabuco'
$data[48,15] = 350
$data[48,16] = 20
$data[48,17] = 'Hortaliza'
$data[49,0] = 9
$data[49,1] = 'Vega Central Mapocho de Santiago'
$data[49,2] = 'Metropolitana'
$data[49,3] = 44433
$data[49,4] = 13
$data[49,5] = 100112005
$data[49,6] = 'Puerro'
$data[49,7] = 'Sin especificar'
$data[49,8] = 'Primera'
$data[49,9] = 142
$data[49,10] = 7000
$data[49,11] = 8000
$data[49,12] = 7500
$data[49,13] = '$/paquete 20 unidades'
$data[49,14] = 'Provincia de Chacabuco'
$data[49,15] = 375
$data[49,16] = 20
$data[49,17] = 'Hortaliza'
$data[50,0] = 9
$data[50,1] = 'Vega Central Mapocho de Santiago'
$data[50,2] = 'Metropolitana'
$data[50,3] = 44203
$data[50,4] = 13
$data[50,5] = 100112005
$data[50,6] = 'Puerro'
$data[50,7] = 'Sin especificar'
$data[50,8] = 'Primera'
$data[50,9] = 50
$data[50,10] = 7000
$data[50,11] = 8000
$data[50,12] = 7400
$data[50,13] = '$/paquete 20 unidades'
$data[50,14] = 'Provincia de Chacabuco'
$data[50,15] = 370
$data[50,16] = 20
$data[50,17] = 'Hortaliza'
$data[51,0] = 9
$data[51,1] = 'Vega Central Mapocho de Santiago'
$data[51,2] = 'Metropolitana'
$data[51,3] = 44707
$data[51,4] = 13
$data[51,5] = 100112005
$data[51,6] = 'Puerro'
$data[51,7] = 'Sin especificar'
$data[51,8] = 'Primera'
$data[51,9] = 70
$data[51,10] = 6000
$data[51,11] = 7000
$data[51,12] = 6571
$data[51,13] = '$/paquete 20 unidades'
$data[51,14] = 'Provincia de Chacabuco'
$data[51,15] = 329
$data[51,16] = 20
$data[51,17] = 'Hortaliza'
$data[52,0] = 9
$data[52,1] = 'Vega Central Mapocho de Santiago'
$data[52,2] = 'Metropolitana'
$data[52,3] = 44328
$data[52,4] = 13
$data[52,5] = 100112005
$data[52,6] = 'Puerro'
$data[52,7] = 'Sin especificar'
$data[52,8] = 'Primera'
$data[52,9] = 160
$data[52,10] = 8000
$data[52,11] = 8000
$data[52,12] = 8000
$data[52,13] = '$/paquete 20 unidades'
$data[52,14] = 'Provincia de Chacabuco'
$data[52,15] = 400
$data[52,16] = 20
$data[52,17] = 'Hortaliza'
$data[53,0] = 9
$data[53,1] = 'Vega Central Mapocho de Santiago'
$data[53,2] = 'Metropolitana'
$data[53,3] = 44419
$data[53,4] = 13
$data[53,5] = 100112005
$data[53,6] = 'Puerro'
$data[53,7] = 'Sin especificar'
$data[53,8] = 'Primera'
$data[53,9] = 160
$data[53,10] = 7000
$data[53,11] = 8000
$data[53,12] = 7500
$data[53,13] = '$/paquete 20 unidades'
$data[53,14] = 'Provincia de Chacabuco'
$data[53,15] = 375
$data[53,16] = 20
$data[53,17] = 'Hortaliza'
$data[54,0] = 9
$data[54,1] = 'Vega Central Mapocho de Santiago'
$data[54,2] = 'Metropolitana'
$data[54,3] = 44545
$data[54,4] = 13
$data[54,5] = 100112005
$data[54,6] = 'Puerro'
$data[54,7] = 'Sin especificar'
$data[54,8] = 'Primera'
$data[54,9] = 160
$data[54,10] = 6000
$data[54,11] = 7000
$data[54,12] = 6500
$data[54,13] = '$/paquete 20 unidades'
$data[54,14] = 'Provincia de Chacabuco'
$data[54,15] = 325
$data[54,16] = 20
$data[54,17] = 'Hortaliza'
$data[55,0] = 9
$data[55,1] = 'Vega Central Mapocho de Santiago'
$data[55,2] = 'Metropolitana'
$data[55,3] = 44229
$data[55,4] = 13
$data[55,5] = 100112005
$data[55,6] = 'Puerro'
$data[55,7] = 'Sin especificar'
$data[55,8] = 'Primera'
$data[55,9] = 50
$data[55,10] = 7000
$data[55,11] = 7000
$data[55,12] = 7000
$data[55,13] = '$/paquete 20 unidades'
$data[55,14] = 'Provincia de Chacabuco'
$data[55,15] = 350
$data[55,16] = 20
$data[55,17] = 'Hortaliza'
$data[56,0] = 9
$data[56,1] = 'Vega Central Mapocho de Santiago'
$data[56,2] = 'Metropolitana'
$data[56,3] = 45035
$data[56,4] = 13
$data[56,5] = 100112005
$data[56,6] = 'Puerro'
$data[56,7] = 'Sin especificar'
$data[56,8] = 'Primera'
$data[56,9] = 70
$data[56,10] = 9000
$data[56,11] = 9000
$data[56,12] = 9000
$data[56,13] = '$/paquete 20 unidades'
$data[56,14] = 'Provincia de Chacabuco'
$data[56,15] = 450
$data[56,16] = 20
$data[56,17] = 'Hortaliza'
$data[57,0] = 9
$data[57,1] = 'Vega Central Mapocho de Santiago'
$data[57,2] = 'Metropolitana'
$data[57,3] = 45063
$data[57,4] = 13
$data[57,5] = 100112005
$data[57,6] = 'Puerro'
$data[57,7] = 'Sin especificar'
$data[57,8] = 'Primera'
$data[57,9] = 70
$data[57,10] = 7000
$data[57,11] = 7000
$data[57,12] = 7000
$data[57,13] = '$/paquete 20 unidades'
$data[57,14] = 'Provincia de Chacabuco'
$data[57,15] = 350
$data[57,16] = 20
$data[57,17] = 'Hortaliza'
$data[58,0] = 9
$data[58,1] = 'Vega Central Mapocho de Santiago'
$data[58,2] = 'Metropolitana'
$data[58,3] = 45030
$data[58,4] = 13
$data[58,5] = 100112005
$data[58,6] = 'Puerro'
$data[58,7] = 'Sin especificar'
$data[58,8] = 'Primera'
$data[58,9] = 70
$data[58,10] = 9000
$data[58,11] = 9000
$data[58,12] = 9000
$data[58,13] = '$/paquete 20 unidades'
$data[58,14] = 'Provincia de Chacabuco'
$data[58,15] = 450
$data[58,16] = 20
$data[58,17] = 'Hortaliza'
$data[59,0] = 9
$data[59,1] = 'Vega Central Mapocho de Santiago'
$data[59,2] = 'Metropolitana'
$data[59,3] = 44252
$data[59,4] = 13
$data[59,5] = 100112005
$data[59,6] = 'Puerro'
$data[59,7] = 'Sin especificar'
$data[59,8] = 'Primera'
$data[59,9] = 160
$data[59,10] = 8000
$data[59,11] = 8000
$data[59,12] = 8000
$data[59,13] = '$/paquete 20 unidades'
$data[59,14] = 'Provincia de Chacabuco'
$data[59,15] = 400
$data[59,16] = 20
$data[59,17] = 'Hortaliza'
$data[60,0] = 9
$data[60,1] = 'Vega Central Mapocho de Santiago'
$data[60,2] = 'Metropolitana'
$data[60,3] = 44321
$data[60,4] = 13
$data[60,5] = 100112005
$data[60,6] = 'Puerro'
$data[60,7] = 'Sin especificar'
$data[60,8] = 'Primera'
$data[60,9] = 250
$data[60,10] = 7000
$data[60,11] = 7000
$data[60,12] = 7000
$data[60,13] = '$/paquete 20 unidades'
$data[60,14] = 'Provincia de Chacabuco'
$data[60,15] = 350
$data[60,16] = 20
$data[60,17] = 'Hortaliza'
$data[61,0] = 9
$data[61,1] = 'Vega Central Mapocho de Santiago'
$data[61,2] = 'Metropolitana'
$data[61,3] = 44267
$data[61,4] = 13
$data[61,5] = 100112005
$data[61,6] = 'Puerro'
$data[61,7] = 'Sin especificar'
$data[61,8] = 'Primera'
$data[61,9] = 160
$data[61,10] = 8000
$data[61,11] = 8000
$data[61,12] = 8000
$data[61,13] = '$/paquete 20 unidades'
$data[61,14] = 'Provincia de Chacabuco'
$data[61,15] = 400
$data[61,16] = 20
$data[61,17] = 'Hortaliza'
$data[62,0] = 9
$data[62,1] = 'Vega Central Mapocho de Santiago'
$data[62,2] = 'Metropolitana'
$data[62,3] = 44302
$data[62,4] = 13
$data[62,5] = 100112005
$data[62,6] = 'Puerro'
$data[62,7] = 'Sin especificar'
$data[62,8] = 'Primera'
$data[62,9] = 160
$data[62,10] = 8000
$data[62,11] = 8000
$data[62,12] = 8000
$data[62,13] = '$/paquete 20 unidades'
$data[62,14] = 'Provincia de Chacabuco'
$data[62,15] = 400
$data[62,16] = 20
$data[62,17] = 'Hortaliza'
$data[63,0] = 9
$data[63,1] = 'Vega Central Mapocho de Santiago'
$data[63,2] = 'Metropolitana'
$data[63,3] = 45042
$data[63,4] = 13
$data[63,5] = 100112005
$data[63,6] = 'Puerro'
$data[63,7] = 'Sin especificar'
$data[63,8] = 'Primera'
$data[63,9] = 70
$data[63,10] = 7000
$data[63,11] = 7000
$data[63,12] = 7000
$data[63,13] = '$/paquete 20 unidades'
$data[63,14] = 'Provincia de Chacabuco'
$data[63,15] = 350
$data[63,16] = 20
$data[63,17] = 'Hortaliza'
$data[64,0] = 9
$data[64,1] = 'Vega Central Mapocho de Santiago'
$data[64,2] = 'Metropolitana'
$data[64,3] = 44244
$data[64,4] = 13
$data[64,5] = 100112005
$data[64,6] = 'Puerro'
$data[64,7] = 'Sin especificar'
$data[64,8] = 'Primera'
$data[64,9] = 70
$data[64,10] = 8000
$data[64,11] = 8000
$data[64,12] = 8000
$data[64,13] = '$/paquete 20 unidades'
$data[64,14] = 'Provincia de Chacabuco'
$data[64,15] = 400
$data[64,16] = 20
$data[64,17] = 'Hortaliza'
$data[65,0] = 9
$data[65,1] = 'Vega Central Mapocho de Santiago'
$data[65,2] = 'Metropolitana'
$data[65,3] = 44162
$data[65,4] = 13
$data[65,5] = 100112005
$data[65,6] = 'Puerro'
$data[65,7] = 'Sin especificar'
$data[65,8] = 'Primera'
$data[65,9] = 50
$data[65,10] = 8000
$data[65,11] = 8000
$data[65,12] = 8000
$data[65,13] = '$/paquete 20 unidades'
$data[65,14] = 'Provincia de Chacabuco'
$data[65,15] = 400
$data[65,16] = 20
$data[65,17] = 'Hortaliza'
$data[66,0] = 9
$data[66,1] = 'Vega Central Mapocho de Santiago'
$data[66,2] = 'Metropolitana'
$data[66,3] = 44855
$data[66,4] = 13
$data[66,5] = 100112005
$data[66,6] = 'Puerro'
$data[66,7] = 'Sin especificar'
$data[66,8] = 'Primera'
$data[66,9] = 65
$data[66,10] = 9000
$data[66,11] = 10000
$data[66,12] = 9538
$data[66,13] = '$/paquete 20 unidades'
$data[66,14] = 'Provincia de Melipilla'
$data[66,15] = 477
$data[66,16] = 20
$data[66,17] = 'Hortaliza'
$data[67,0] = 9
$data[67,1] = 'Vega Central Mapocho de Santiago'
$data[67,2] = 'Metropolitana'
$data[67,3] = 44855
$data[67,4] = 13
$data[67,5] = 100112005
$data[67,6] = 'Puerro'
$data[67,7] = 'Sin especificar'
$data[67,8] = 'Primera'
$data[67,9] = 80
$data[67,10] = 10000
$data[67,11] = 10000
$data[67,12] = 10000
$data[67,13] = '$/paquete 20 unidades'
$data[67,14] = 'Provincia de Santiago'
$data[67,15] = 500
$data[67,16] = 20
$data[67,17] = 'Hortaliza'
$data[68,0] = 9
$data[68,1] = 'Vega Central Mapocho de Santiago'
$data[68,2] = 'Metropolitana'
$data[68,3] = 44329
$data[68,4] = 13
$data[68,5] = 100112005
$data[68,6] = 'Puerro'
$data[68,7] = 'Sin especificar'
$data[68,8] = 'Primera'
$data[68,9] = 160
$data[68,10] = 8000
$data[68,11] = 8000
$data[68,12] = 8000
$data[68,13] = '$/paquete 20 unidades'
$data[68,14] = 'Provincia de Chacabuco'
$data[68,15] = 400
$data[68,16] = 20
$data[68,17] = 'Hortaliza'
$data[69,0] = 9
$data[69,1] = 'Vega Central Mapocho de Santiago'
$data[69,2] = 'Metropolitana'
$data[69,3] = 44405
$data[69,4] = 13
$data[69,5] = 100112005
$data[69,6] = 'Puerro'
$data[69,7] = 'Sin especificar'
$data[69,8] = 'Primera'
$data[69,9] = 160
$data[69,10] = 7000
$data[69,11] = 8000
$data[69,12] = 7500
$data[69,13] = '$/paquete 20 unidades'
$data[69,14] = 'Provincia de Chacabuco'
$data[69,15] = 375
$data[69,16] = 20
$data[69,17] = 'Hortaliza'
$data[70,0] = 9
$data[70,1] = 'Vega Central Mapocho de Santiago'
$data[70,2] = 'Metropolitana'
$data[70,3] = 44273
$data[70,4] = 13
$data[70,5] = 100112005
$data[70,6] = 'Puerro'
$data[70,7] = 'Sin especificar'
$data[70,8] = 'Primera'
$data[70,9] = 70
$data[70,10] = 8000
$data[70,11] = 8000
$data[70,12] = 8000
$data[70,13] = '$/paquete 20 unidades'
$data[70,14] = 'Provincia de Chacabuco'
$data[70,15] = 400
$data[70,16] = 20
$data[70,17] = 'Hortaliza'
$data[71,0] = 9
$data[71,1] = 'Vega Central Mapocho de Santiago'
$data[71,2] = 'Metropolitana'
$data[71,3] = 44358
$data[71,4] = 13
$data[71,5] = 100112005
$data[71,6] = 'Puerro'
$data[71,7] = 'Sin especificar'
$data[71,8] = 'Primera'
$data[71,9] = 160
$data[71,10] = 7500
$data[71,11] = 8000
$data[71,12] = 7750
$data[71,13] = '$/paquete 20 unidades'
$data[71,14] = 'Provincia de Chacabuco'
$data[71,15] = 388
$data[71,16] = 20
$data[71,17] = 'Hortaliza'
$data[72,0] = 9
$data[72,1] = 'Vega Central Mapocho de Santiago'
$data[72,2] = 'Metropolitana'
$data[72,3] = 44412
$data[72,4] = 13
$data[72,5] = 100112005
$data[72,6] = 'Puerro'
$data[72,7] = 'Sin especificar'
$data[72,8] = 'Primera'
$data[72,9] = 160
$data[72,10] = 7500
$data[72,11] = 8000
$data[72,12] = 7750
$data[72,13] = '$/paquete 20 unidades'
$data[72,14] = 'Provincia de Chacabuco'
$data[72,15] = 388
$data[72,16] = 20
$data[72,17] = 'Hortaliza'
$data[73,0] = 9
$data[73,1] = 'Vega Central Mapocho de Santiago'
$data[73,2] = 'Metropolitana'
$data[73,3] = 44832
$data[73,4] = 13
$data[73,5] = 100112005
$data[73,6] = 'Puerro'
$data[73,7] = 'Sin especificar'
$data[73,8] = 'Segunda'
$data[73,9] = 30
$data[73,10] = 12000
$data[73,11] = 12000
$data[73,12] = 12000
$data[73,13] = '$/paquete 20 unidades'
$data[73,14] = 'Provincia de Melipilla'
$data[73,15] = 600
$data[73,16] = 20
$data[73,17] = 'Hortaliza'
$data[74,0] = 9
$data[74,1] = 'Vega Central Mapocho de Santiago'
$data[74,2] = 'Metropolitana'
$data[74,3] = 44573
$data[74,4] = 13
$data[74,5] = 100112005
$data[74,6] = 'Puerro'
$data[74,7] = 'Sin especificar'
$data[74,8] = 'Primera'
$data[74,9] = 106
$data[74,10] = 7000
$data[74,11] = 8000
$data[74,12] = 7500
$data[74,13] = '$/paquete 20 unidades'
$data[74,14] = 'Provincia de Chacabuco'
$data[74,15] = 375
$data[74,16] = 20
$data[74,17] = 'Hortaliza'
$data[75,0] = 9
$data[75,1] = 'Vega Central Mapocho de Santiago'
$data[75,2] = 'Metropolitana'
$data[75,3] = 44776
$data[75,4] = 13
$data[75,5] = 100112005
$data[75,6] = 'Puerro'
$data[75,7] = 'Sin especificar'
$data[75,8] = 'Primera'
$data[75,9] = 160
$data[75,10] = 7000
$data[75,11] = 8000
$data[75,12] = 7500
$data[75,13] = '$/paquete 20 unidades'
$data[75,14] = 'Provincia de Chacabuco'
$data[75,15] = 375
$data[75,16] = 20
$data[75,17] = 'Hortaliza'
$data[76,0] = 9
$data[76,1] = 'Vega Central Mapocho de Santiago'
$data[76,2] = 'Metropolitana'
$data[76,3] = 44215
$data[76,4] = 13
$data[76,5] = 100112005
$data[76,6] = 'Puerro'
$data[76,7] = 'Sin especificar'
$data[76,8] = 'Primera'
$data[76,9] = 80
$data[76,10] = 7000
$data[76,11] = 7000
$data[76,12] = 7000
$data[76,13] = '$/paquete 20 unidades'
$data[76,14] = 'Provincia de Chacabuco'
$data[76,15] = 350
$data[76,16] = 20
$data[76,17] = 'Hortaliza'
$data[77,0] = 9
$data[77,1] = 'Vega Central Mapocho de Santiago'
$data[77,2] = 'Metropolitana'
$data[77,3] = 44426
$data[77,4] = 13
$data[77,5] = 100112005
$data[77,6] = 'Puerro'
$data[77,7] = 'Sin especificar'
$data[77,8] = 'Primera'
$data[77,9] = 97
$data[77,10] = 7000
$data[77,11] = 8000
$data[77,12] = 7505
$data[77,13] = '$/paquete 20 unidades'
$data[77,14] = 'Provincia de Chacabuco'
$data[77,15] = 375
$data[77,16] = 20
$data[77,17] = 'Hortaliza'
$data[78,0] = 9
$data[78,1] = 'Vega Central Mapocho de Santiago'
$data[78,2] = 'Metropolitana'
$data[78,3] = 44160
$data[78,4] = 13
$data[78,5] = 100112005
$data[78,6] = 'Puerro'
$data[78,7] = 'Sin especificar'
$data[78,8] = 'Primera'
$data[78,9] = 50
$data[78,10] = 7000
$data[78,11] = 8000
$data[78,12] = 7600
$data[78,13] = '$/paquete 20 unidades'
$data[78,14] = 'Provincia de Chacabuco'
$data[78,15] = 380
$data[78,16] = 20
$data[78,17] = 'Hortaliza'
$data[79,0] = 9
$data[79,1] = 'Vega Central Mapocho de Santiago'
$data[79,2] = 'Metropolitana'
$data[79,3] = 45028
$data[79,4] = 13
$data[79,5] = 100112005
$data[79,6] = 'Puerro'
$data[79,7] = 'Sin especificar'
$data[79,8] = 'Primera'
$data[79,9] = 160
$data[79,10] = 8000
$data[79,11] = 8000
$data[79,12] = 8000
$data[79,13] = '$/paquete 20 unidades'
$data[79,14] = 'Provincia de Chacabuco'
$data[79,15] = 400
$data[79,16] = 20
$data[79,17] = 'Hortaliza'
$data[80,0] = 9
$data[80,1] = 'Vega Central Mapocho de Santiago'
$data[80,2] = 'Metropolitana'
$data[80,3] = 44442
$data[80,4] = 13
$data[80,5] = 100112005
$data[80,6] = 'Puerro'
$data[80,7] = 'Sin especificar'
$data[80,8] = 'Primera'
$data[80,9] = 180
$data[80,10] = 7000
$data[80,11] = 8000
$data[80,12] = 7500
$data[80,13] = '$/paquete 20 unidades'
$data[80,14] = 'Provincia de Chacabuco'
$data[80,15] = 375
$data[80,16] = 20
$data[80,17] = 'Hortaliza'
$data[81,0] = 9
$data[81,1] = 'Vega Central Mapocho de Santiago'
$data[81,2] = 'Metropolitana'
$data[81,3] = 44685
$data[81,4] = 13
$data[81,5] = 100112005
$data[81,6] = 'Puerro'
$data[81,7] = 'Sin especificar'
$data[81,8] = 'Primera'
$data[81,9] = 160
$data[81,10] = 8000
$data[81,11] = 8000
$data[81,12] = 8000
$data[81,13] = '$/paquete 20 unidades'
$data[81,14] = 'Provincia de Chacabuco'
$data[81,15] = 400
$data[81,16] = 20
$data[81,17] = 'Hortaliza'
$data[82,0] = 9
$data[82,1] = 'Vega Central Mapocho de Santiago'
$data[82,2] = 'Metropolitana'
$data[82,3] = 44363
$data[82,4] = 13
$data[82,5] = 100112005
$data[82,6] = 'Puerro'
$data[82,7] = 'Sin especificar'
$data[82,8] = 'Primera'
$data[82,9] = 160
$data[82,10] = 8000
$data[82,11] = 8000
$data[82,12] = 8000
$data[82,13] = '$/paquete 20 unidades'
$data[82,14] = 'Provincia de Chacabuco'
$data[82,15] = 400
$data[82,16] = 20
$data[82,17] = 'Hortaliza'
$data[83,0] = 9
$data[83,1] = 'Vega Central Mapocho de Santiago'
$data[83,2] = 'Metropolitana'
$data[83,3] = 45084
$data[83,4] = 13
$data[83,5] = 100112005
$data[83,6] = 'Puerro'
$data[83,7] = 'Sin especificar'
$data[83,8] = 'Primera'
$data[83,9] = 160
$data[83,10] = 7000
$data[83,11] = 7000
$data[83,12] = 7000
$data[83,13] = '$/paquete 20 unidades'
$data[83,14] = 'Provincia de Chacabuco'
$data[83,15] = 350
$data[83,16] = 20
$data[83,17] = 'Hortaliza'
$data[84,0] = 9
$data[84,1] = 'Vega Central Mapocho de Santiago'
$data[84,2] = 'Metropolitana'
$data[84,3] = 44281
$data[84,4] = 13
$data[84,5] = 100112005
$data[84,6] = 'Puerro'
$data[84,7] = 'Sin especificar'
$data[84,8] = 'Primera'
$data[84,9] = 250
$data[84,10] = 8000
$data[84,11] = 8000
$data[84,12] = 8000
$data[84,13] = '$/paquete 20 unidades'
$data[84,14] = 'Provincia de Chacabuco'
$data[84,15] = 400
$data[84,16] = 20
$data[84,17] = 'Hortaliza'
$data[85,0] = 9
$data[85,1] = 'Vega Central Mapocho de Santiago'
$data[85,2] = 'Metropolitana'
$data[85,3] = 44259
$data[85,4] = 13
$data[85,5] = 100112005
$data[85,6] = 'Puerro'
$data[85,7] = 'Sin especificar'
$data[85,8] = 'Primera'
$data[85,9] = 160
$data[85,10] = 8000
$data[85,11] = 8000
$data[85,12] = 8000
$data[85,13] = '$/paquete 20 unidades'
$data[85,14] = 'Provincia de Chacabuco'
$data[85,15] = 400
$data[85,16] = 20
$data[85,17] = 'Hortaliza'
$data[86,0] = 9
$data[86,1] = 'Vega Central Mapocho de Santiago'
$data[86,2] = 'Metropolitana'
$data[86,3] = 44699
$data[86,4] = 13
$data[86,5] = 100112005
$data[86,6] = 'Puerro'
$data[86,7] = 'Sin especificar'
$data[86,8] = 'Primera'
$data[86,9] = 160
$data[86,10] = 7000
$data[86,11] = 8000
$data[86,12] = 7500
$data[86,13] = '$/paquete 20 unidades'
$data[86,14] = 'Provincia de Chacabuco'
$data[86,15] = 375
$data[86,16] = 20
$data[86,17] = 'Hortaliza'
$data[87,0] = 9
$data[87,1] = 'Vega Central Mapocho de Santiago'
$data[87,2] = 'Metropolitana'
$data[87,3] = 44484
$data[87,4] = 13
$data[87,5] = 100112005
$data[87,6] = 'Puerro'
$data[87,7] = 'Sin especificar'
$data[87,8] = 'Primera'
$data[87,9] = 160
$data[87,10] = 7000
$data[87,11] = 8000
$data[87,12] = 7500
$data[87,13] = '$/paquete 20 unidades'
$data[87,14] = 'Provincia de Chacabuco'
$data[87,15] = 375
$data[87,16] = 20
$data[87,17] = 'Hortaliza'
$data[88,0] = 9
$data[88,1] = 'Vega Central Mapocho de Santiago'
$data[88,2] = 'Metropolitana'
$data[88,3] = 44552
$data[88,4] = 13
$data[88,5] = 100112005
$data[88,6] = 'Puerro'
$data[88,7] = 'Sin especificar'
$data[88,8] = 'Primera'
$data[88,9] = 106
$data[88,10] = 7000
$data[88,11] = 8000
$data[88,12] = 7500
$data[88,13] = '$/paquete 20 unidades'
$data[88,14] = 'Provincia de Chacabuco'
$data[88,15] = 375
$data[88,16] = 20
$data[88,17] = 'Hortaliza'
$data[89,0] = 9
$data[89,1] = 'Vega Central Mapocho de Santiago'
$data[89,2] = 'Metropolitana'
$data[89,3] = 44272
$data[89,4] = 13
$data[89,5] = 100112005
$data[89,6] = 'Puerro'
$data[89,7] = 'Sin especificar'
$data[89,8] = 'Primera'
$data[89,9] = 160
$data[89,10] = 8000
$data[89,11] = 8000
$data[89,12] = 8000
$data[89,13] = '$/paquete 20 unidades'
$data[89,14] = 'Provincia de Chacabuco'
$data[89,15] = 400
$data[89,16] = 20
$data[89,17] = 'Hortaliza'
$data[90,0] = 9
$data[90,1] = 'Vega Central Mapocho de Santiago'
$data[90,2] = 'Metropolitana'
$data[90,3] = 44623
$data[90,4] = 13
$data[90,5] = 100112005
$data[90,6] = 'Puerro'
$data[90,7] = 'Sin especificar'
$data[90,8] = 'Primera'
$data[90,9] = 106
$data[90,10] = 7000
$data[90,11] = 8000
$data[90,12] = 7500
$data[90,13] = '$/paquete 20 unidades'
$data[90,14] = 'Provincia de Chacabuco'
$data[90,15] = 375
$data[90,16] = 20
$data[90,17] = 'Hortaliza'
$data[91,0] = 9
$data[91,1] = 'Vega Central Mapocho de Santiago'
$data[91,2] = 'Metropolitana'
$data[91,3] = 44664
$data[91,4] = 13
$data[91,5] = 100112005
$data[91,6] = 'Puerro'
$data[91,7] = 'Sin especificar'
$data[91,8] = 'Primera'
$data[91,9] = 106
$data[91,10] = 8000
$data[91,11] = 8000
$data[91,12] = 8000
$data[91,13] = '$/paquete 20 unidades'
$data[91,14] = 'Provincia de Chacabuco'
$data[91,15] = 400
$data[91,16] = 20
$data[91,17] = 'Hortaliza'
$data[92,0] = 9
$data[92,1] = 'Vega Central Mapocho de Santiago'
$data[92,2] = 'Metropolitana'
$data[92,3] = 44510
$data[92,4] = 13
$data[92,5] = 100112005
$data[92,6] = 'Puerro'
$data[92,7] = 'Sin especificar'
$data[92,8] = 'Primera'
$data[92,9] = 160
$data[92,10] = 6000
$data[92,11] = 7000
$data[92,12] = 6500
$data[92,13] = '$/paquete 20 unidades'
$data[92,14] = 'Provincia de Chacabuco'
$data[92,15] = 325
$data[92,16] = 20
$data[92,17] = 'Hortaliza'
$data[93,0] = 9
$data[93,1] = 'Vega Central Mapocho de Santiago'
$data[93,2] = 'Metropolitana'
$data[93,3] = 44636
$data[93,4] = 13
$data[93,5] = 100112005
$data[93,6] = 'Puerro'
$data[93,7] = 'Sin especificar'
$data[93,8] = 'Primera'
$data[93,9] = 97
$data[93,10] = 7000
$data[93,11] = 7000
$data[93,12] = 7000
$data[93,13] = '$/paquete 20 unidades'
$data[93,14] = 'Provincia de Chacabuco'
$data[93,15] = 350
$data[93,16] = 20
$data[93,17] = 'Hortaliza'
$data[94,0] = 9
$data[94,1] = 'Vega Central Mapocho de Santiago'
$data[94,2] = 'Metropolitana'
$data[94,3] = 45091
$data[94,4] = 13
$data[94,5] = 100112005
$data[94,6] = 'Puerro'
$data[94,7] = 'Sin especificar'
$data[94,8] = 'Primera'
$data[94,9] = 160
$data[94,10] = 7000
$data[94,11] = 8000
$data[94,12] = 7500
$data[94,13] = '$/paquete 20 unidades'
$data[94,14] = 'Provincia de Chacabuco'
$data[94,15] = 375
$data[94,16] = 20
$data[94,17] = 'Hortaliza'
$data[95,0] = 9
$data[95,1] = 'Vega Central Mapocho de Santiago'
$data[95,2] = 'Metropolitana'
$data[95,3] = 44214
$data[95,4] = 13
$data[95,5] = 100112005
$data[95,6] = 'Puerro'
$data[95,7] = 'Sin especificar'
$data[95,8] = 'Primera'
$data[95,9] = 50
$data[95,10] = 8000
$data[95,11] = 8000
$data[95,12] = 8000
$data[95,13] = '$/paquete 20 unidades'
$data[95,14] = 'Provincia de Chacabuco'
$data[95,15] = 400
$data[95,16] = 20
$data[95,17] = 'Hortaliza'
$data[96,0] = 9
$data[96,1] = 'Vega Central Mapocho de Santiago'
$data[96,2] = 'Metropolitana'
$data[96,3] = 44306
$data[96,4] = 13
$data[96,5] = 100112005
$data[96,6] = 'Puerro'
$data[96,7] = 'Sin especificar'
$data[96,8] = 'Primera'
$data[96,9] = 160
$data[96,10] = 8000
$data[96,11] = 8000
$data[96,12] = 8000
$data[96,13] = '$/paquete 20 unidades'
$data[96,14] = 'Provincia de Chacabuco'
$data[96,15] = 400
$data[96,16] = 20
$data[96,17] = 'Hortaliza'
$data[97,0] = 9
$data[97,1] = 'Vega Central Mapocho de Santiago'
$data[97,2] = 'Metropolitana'
$data[97,3] = 44482
$data[97,4] = 13
$data[97,5] = 100112005
$data[97,6] = 'Puerro'
$data[97,7] = 'Sin especificar'
$data[97,8] = 'Primera'
$data[97,9] = 160
$data[97,10] = 7000
$data[97,11] = 8000
$data[97,12] = 7500
$data[97,13] = '$/paquete 20 unidades'
$data[97,14] = 'Provincia de Chacabuco'
$data[97,15] = 375
$data[97,16] = 20
$data[97,17] = 'Hortaliza'
$data[98,0] = 9
$data[98,1] = 'Vega Central Mapocho de Santiago'
$data[98,2] = 'Metropolitana'
$data[98,3] = 44497
$data[98,4] = 13
$data[98,5] = 100112005
$data[98,6] = 'Puerro'
$data[98,7] = 'Sin especificar'
$data[98,8] = 'Primera'
$data[98,9] = 180
$data[98,10] = 6000
$data[98,11] = 7000
$data[98,12] = 6556
$data[98,13] = '$/paquete 20 unidades'
$data[98,14] = 'Provincia de Chacabuco'
$data[98,15] = 328
$data[98,16] = 20
$data[98,17] = 'Hortaliza'
$data[99,0] = 9
$data[99,1] = 'Vega Central Mapocho de Santiago'
$data[99,2] = 'Metropolitana'
$data[99,3] = 44762
$data[99,4] = 13
$data[99,5] = 100112005
$data[99,6] = 'Puerro'
$data[99,7] = 'Sin especificar'
$data[99,8] = 'Primera'
$data[99,9] = 160
$data[99,10] = 7000
$data[99,11] = 8000
$data[99,12] = 7500
$data[99,13] = '$/paquete 20 unidades'
$data[99,14] = 'Provincia de Chacabuco'
$data[99,15] = 375
$data[99,16] = 20
$data[99,17] = 'Hortaliza'
$data[100,0] = 9
$data[100,1] = 'Vega Central Mapocho de Santiago'
$data[100,2] = 'Metropolitana'
$data[100,3] = 44189
$data[100,4] = 13
$data[100,5] = 100112005
$data[100,6] = 'Puerro'
$data[100,7] = 'Sin especificar'
$data[100,8] = 'Primera'
$data[100,9] = 50
$data[100,10] = 8000
$data[100,11] = 8000
$data[100,12] = 8000
$data[100,13] = '$/paquete 20 unidades'
$data[100,14] = 'Provincia de Chacabuco'
$data[100,15] = 400
$data[100,16] = 20
$data[100,17] = 'Hortaliza'
$data[101,0] = 9
$data[101,1] = 'Vega Central Mapocho de Santiago'
$data[101,2] = 'Metropolitana'
$data[101,3] = 44398
$data[101,4] = 13
$data[101,5] = 100112005
$data[101,6] = 'Puerro'
$data[101,7] = 'Sin especificar'
$data[101,8] = 'Primera'
$data[101,9] = 70
$data[101,10] = 7500
$data[101,11] = 8000
$data[101,12] = 7750
$data[101,13] = '$/paquete 20 unidades'
$data[101,14] = 'Provincia de Chacabuco'
$data[101,15] = 388
$data[101,16] = 20
$data[101,17] = 'Hortaliza'
$data[102,0] = 9
$data[102,1] = 'Vega Central Mapocho de Santiago'
$data[102,2] = 'Metropolitana'
$data[102,3] = 44218
$data[102,4] = 13
$data[102,5] = 100112005
$data[102,6] = 'Puerro'
$data[102,7] = 'Sin especificar'
$data[102,8] = 'Primera'
$data[102,9] = 80
$data[102,10] = 6000
$data[102,11] = 7000
$data[102,12] = 6625
$data[102,13] = '$/paquete 20 unidades'
$data[102,14] = 'Provincia de Chacabuco'
$data[102,15] = 331
$data[102,16] = 20
$data[102,17] = 'Hortaliza'
$data[103,0] = 9
$data[103,1] = 'Vega Central Mapocho de Santiago'
$data[103,2] = 'Metropolitana'
$data[103,3] = 44435
$data[103,4] = 13
$data[103,5] = 100112005
$data[103,6] = 'Puerro'
$data[103,7] = 'Sin especificar'
$data[103,8] = 'Primera'
$data[103,9] = 302
$data[103,10] = 7000
$data[103,11] = 8000
$data[103,12] = 7500
$data[103,13] = '$/paquete 20 unidades'
$data[103,14] = 'Provincia de Chacabuco'
$data[103,15] = 375
$data[103,16] = 20
$data[103,17] = 'Hortaliza'
$data[104,0] = 9
$data[104,1] = 'Vega Central Mapocho de Santiago'
$data[104,2] = 'Metropolitana'
$data[104,3] = 44825
$data[104,4] = 13
$data[104,5] = 100112005
$data[104,6] = 'Puerro'
$data[104,7] = 'Sin especificar'
$data[104,8] = 'Primera'
$data[104,9] = 70
$data[104,10] = 12000
$data[104,11] = 12000
$data[104,12] = 12000
$data[104,13] = '$/paquete 20 unidades'
$data[104,14] = 'Provincia de Chacabuco'
$data[104,15] = 600
$data[104,16] = 20
$data[104,17] = 'Hortaliza'
$data[105,0] = 9
$data[105,1] = 'Vega Central Mapocho de Santiago'
$data[105,2] = 'Metropolitana'
$data[105,3] = 44454
$data[105,4] = 13
$data[105,5] = 100112005
$data[105,6] = 'Puerro'
$data[105,7] = 'Sin especificar'
$data[105,8] = 'Primera'
$data[105,9] = 160
$data[105,10] = 7000
$data[105,11] = 8000
$data[105,12] = 7500
$data[105,13] = '$/paquete 20 unidades'
$data[105,14] = 'Provincia de Chacabuco'
$data[105,15] = 375
$data[105,16] = 20
$data[105,17] = 'Hortaliza'
$data[106,0] = 9
$data[106,1] = 'Vega Central Mapocho de Santiago'
$data[106,2] = 'Metropolitana'
$data[106,3] = 44790
$data[106,4] = 13
$data[106,5] = 100112005
$data[106,6] = 'Puerro'
$data[106,7] = 'Sin especificar'
$data[106,8] = 'Primera'
$data[106,9] = 160
$data[106,10] = 7000
$data[106,11] = 7000
$data[106,12] = 7000
$data[106,13] = '$/paquete 20 unidades'
$data[106,14] = 'Provincia de Chacabuco'
$data[106,15] = 350
$data[106,16] = 20
$data[106,17] = 'Hortaliza'
$data[107,0] = 9
$data[107,1] = 'Vega Central Mapocho de Santiago'
$data[107,2] = 'Metropolitana'
$data[107,3] = 44720
$data[107,4] = 13
$data[107,5] = 100112005
$data[107,6] = 'Puerro'
$data[107,7] = 'Sin especificar'
$data[107,8] = 'Primera'
$data[107,9] = 160
$data[107,10] = 7000
$data[107,11] = 8000
$data[107,12] = 7500
$data[107,13] = '$/paquete 20 unidades'
$data[107,14] = 'Provincia de Chacabuco'
$data[107,15] = 375
$data[107,16] = 20
$data[107,17] = 'Hortaliza'
$data[108,0] = 9
$data[108,1] = 'Vega Central Mapocho de Santiago'
$data[108,2] = 'Metropolitana'
$data[108,3] = 44615
$data[108,4] = 13
$data[108,5] = 100112005
$data[108,6] = 'Puerro'
$data[108,7] = 'Sin especificar'
$data[108,8] = 'Primera'
$data[108,9] = 79
$data[108,10] = 7000
$data[108,11] = 7000
$data[108,12] = 7000
$data[108,13] = '$/paquete 20 unidades'
$data[108,14] = 'Provincia de Chacabuco'
$data[108,15] = 350
$data[108,16] = 20
$data[108,17] = 'Hortaliza'
$data[109,0] = 9
$data[109,1] = 'Vega Central Mapocho de Santiago'
$data[109,2] = 'Metropolitana'
$data[109,3] = 44314
$data[109,4] = 13
$data[109,5] = 100112005
$data[109,6] = 'Puerro'
$data[109,7] = 'Sin especificar'
$data[109,8] = 'Primera'
$data[109,9] = 160
$data[109,10] = 8000
$data[109,11] = 8000
$data[109,12] = 8000
$data[109,13] = '$/paquete 20 unidades'
$data[109,14] = 'Provincia de Chacabuco'
$data[109,15] = 400
$data[109,16] = 20
$data[109,17] = 'Hortaliza'
$data[110,0] = 9
$data[110,1] = 'Vega Central Mapocho de Santiago'
$data[110,2] = 'Metropolitana'
$data[110,3] = 44370
$data[110,4] = 13
$data[110,5] = 100112005
$data[110,6] = 'Puerro'
$data[110,7] = 'Sin especificar'
$data[110,8] = 'Primera'
$data[110,9] = 160
$data[110,10] = 7500
$data[110,11] = 8000
$data[110,12] = 7750
$data[110,13] = '$/paquete 20 unidades'
$data[110,14] = 'Provincia de Chacabuco'
$data[110,15] = 388
$data[110,16] = 20
$data[110,17] = 'Hortaliza'
$data[111,0] = 9
$data[111,1] = 'Vega Central Mapocho de Santiago'
$data[111,2] = 'Metropolitana'
$data[111,3] = 44342
$data[111,4] = 13
$data[111,5] = 100112005
$data[111,6] = 'Puerro'
$data[111,7] = 'Sin especificar'
$data[111,8] = 'Primera'
$data[111,9] = 160
$data[111,10] = 8000
$data[111,11] = 8000
$data[111,12] = 8000
$data[111,13] = '$/paquete 20 unidades'
$data[111,14] = 'Provincia de Chacabuco'
$data[111,15] = 400
$data[111,16] = 20
$data[111,17] = 'Hortaliza'
$data[112,0] = 9
$data[112,1] = 'Vega Central Mapocho de Santiago'
$data[112,2] = 'Metropolitana'
$data[112,3] = 44489
$data[112,4] = 13
$data[112,5] = 100112005
$data[112,6] = 'Puerro'
$data[112,7] = 'Sin especificar'
$data[112,8] = 'Primera'
$data[112,9] = 160
$data[112,10] = 7000
$data[112,11] = 8000
$data[112,12] = 7500
$data[112,13] = '$/paquete 20 unidades'
$data[112,14] = 'Provincia de Chacabuco'
$data[112,15] = 375
$data[112,16] = 20
$data[112,17] = 'Hortaliza'
$data[113,0] = 9
$data[113,1] = 'Vega Central Mapocho de Santiago'
$data[113,2] = 'Metropolitana'
$data[113,3] = 44503
$data[113,4] = 13
$data[113,5] = 100112005
$data[113,6] = 'Puerro'
$data[113,7] = 'Sin especificar'
$data[113,8] = 'Primera'
$data[113,9] = 97
$data[113,10] = 7000
$data[113,11] = 8000
$data[113,12] = 7505
$data[113,13] = '$/paquete 20 unidades'
$data[113,14] = 'Provincia de Chacabuco'
$data[113,15] = 375
$data[113,16] = 20
$data[113,17] = 'Hortaliza'
$data[114,0] = 9
$data[114,1] = 'Vega Central Mapocho de Santiago'
$data[114,2] = 'Metropolitana'
$data[114,3] = 45015
$data[114,4] = 13
$data[114,5] = 100112005
$data[114,6] = 'Puerro'
$data[114,7] = 'Sin especificar'
$data[114,8] = 'Primera'
$data[114,9] = 70
$data[114,10] = 7000
$data[114,11] = 8000
$data[114,12] = 7500
$data[114,13] = '$/paquete 20 unidades'
$data[114,14] = 'Provincia de Chacabuco'
$data[114,15] = 375
$data[114,16] = 20
$data[114,17] = 'Hortaliza'
$data[115,0] = 9
$data[115,1] = 'Vega Central Mapocho de Santiago'
$data[115,2] = 'Metropolitana'
$data[115,3] = 44692
$data[115,4] = 13
$data[115,5] = 100112005
$data[115,6] = 'Puerro'
$data[115,7] = 'Sin especificar'
$data[115,8] = 'Primera'
$data[115,9] = 124
$data[115,10] = 7000
$data[115,11] = 8000
$data[115,12] = 7500
$data[115,13] = '$/paquete 20 unidades'
$data[115,14] = 'Provincia de Chacabuco'
$data[115,15] = 375
$data[115,16] = 20
$data[115,17] = 'Hortaliza'
$data[116,0] = 9
$data[116,1] = 'Vega Central Mapocho de Santiago'
$data[116,2] = 'Metropolitana'
$data[116,3] = 44166
$data[116,4] = 13
$data[116,5] = 100112005
$data[116,6] = 'Puerro'
$data[116,7] = 'Sin especificar'
$data[116,8] = 'Primera'
$data[116,9] = 50
$data[116,10] = 8000
$data[116,11] = 8000
$data[116,12] = 8000
$data[116,13] = '$/paquete 20 unidades'
$data[116,14] = 'Provincia de Chacabuco'
$data[116,15] = 400
$data[116,16] = 20
$data[116,17] = 'Hortaliza'
$data[117,0] = 9
$data[117,1] = 'Vega Central Mapocho de Santiago'
$data[117,2] = 'Metropolitana'
$data[117,3] = 45077
$data[117,4] = 13
$data[117,5] = 100112005
$data[117,6] = 'Puerro'
$data[117,7] = 'Sin especificar'
$data[117,8] = 'Primera'
$data[117,9] = 70
$data[117,10] = 8000
$data[117,11] = 8000
$data[117,12] = 8000
$data[117,13] = '$/paquete 20 unidades'
$data[117,14] = 'Provincia de Chacabuco'
$data[117,15] = 400
$data[117,16] = 20
$data[117,17] = 'Hortaliza'
$data[118,0] = 9
$data[118,1] = 'Vega Central Mapocho de Santiago'
$data[118,2] = 'Metropolitana'
$data[118,3] = 44491
$data[118,4] = 13
$data[118,5] = 100112005
$data[118,6] = 'Puerro'
$data[118,7] = 'Sin especificar'
$data[118,8] = 'Primera'
$data[118,9] = 160
$data[118,10] = 7000
$data[118,11] = 8000
$data[118,12] = 7500
$data[118,13] = '$/paquete 20 unidades'
$data[118,14] = 'Provincia de Chacabuco'
$data[118,15] = 375
$data[118,16] = 20
$data[118,17] = 'Hortaliza'
$data[119,0] = 9
$data[119,1] = 'Vega Central Mapocho de Santiago'
$data[119,2] = 'Metropolitana'
$data[119,3] = 44407
$data[119,4] = 13
$data[119,5] = 100112005
$data[119,6] = 'Puerro'
$data[119,7] = 'Sin especificar'
$data[119,8] = 'Primera'
$data[119,9] = 160
$data[119,10] = 7000
$data[119,11] = 8000
$data[119,12] = 7500
$data[119,13] = '$/paquete 20 unidades'
$data[119,14] = 'Provincia de Chacabuco'
$data[119,15] = 375
$data[119,16] = 20
$data[119,17] = 'Hortaliza'
$data[120,0] = 9
$data[120,1] = 'Vega Central Mapocho de Santiago'
$data[120,2] = 'Metropolitana'
$data[120,3] = 44860
$data[120,4] = 13
$data[120,5] = 100112005
$data[120,6] = 'Puerro'
$data[120,7] = 'Sin especificar'
$data[120,8] = 'Primera'
$data[120,9] = 70
$data[120,10] = 9000
$data[120,11] = 10000
$data[120,12] = 9571
$data[120,13] = '$/paquete 20 unidades'
$data[120,14] = 'Provincia de Chacabuco'
$data[120,15] = 479
$data[120,16] = 20
$data[120,17] = 'Hortaliza'
$data[121,0] = 9
$data[121,1] = 'Vega Central Mapocho de Santiago'
$data[121,2] = 'Metropolitana'
$data[121,3] = 44447
$data[121,4] = 13
$data[121,5] = 100112005
$data[121,6] = 'Puerro'
$data[121,7] = 'Sin especificar'
$data[121,8] = 'Primera'
$data[121,9] = 106
$data[121,10] = 7000
$data[121,11] = 8000
$data[121,12] = 7500
$data[121,13] = '$/paquete 20 unidades'
$data[121,14] = 'Provincia de Chacabuco'
$data[121,15] = 375
$data[121,16] = 20
$data[121,17] = 'Hortaliza'
$data[122,0] = 9
$data[122,1] = 'Vega Central Mapocho de Santiago'
$data[122,2] = 'Metropolitana'
$data[122,3] = 44292
$data[122,4] = 13
$data[122,5] = 100112005
$data[122,6] = 'Puerro'
$data[122,7] = 'Sin especificar'
$data[122,8] = 'Primera'
$data[122,9] = 70
$data[122,10] = 8000
$data[122,11] = 8000
$data[122,12] = 8000
$data[122,13] = '$/paquete 20 unidades'
$data[122,14] = 'Provincia de Chacabuco'
$data[122,15] = 400
$data[122,16] = 20
$data[122,17] = 'Hortaliza'
$data[123,0] = 9
$data[123,1] = 'Vega Central Mapocho de Santiago'
$data[123,2] = 'Metropolitana'
$data[123,3] = 44295
$data[123,4] = 13
$data[123,5] = 100112005
$data[123,6] = 'Puerro'
$data[123,7] = 'Sin especificar'
$data[123,8] = 'Primera'
$data[123,9] = 70
$data[123,10] = 8000
$data[123,11] = 8000
$data[123,12] = 8000
$data[123,13] = '$/paquete 20 unidades'
$data[123,14] = 'Provincia de Chacabuco'
$data[123,15] = 400
$data[123,16] = 20
$data[123,17] = 'Hortaliza'
$data[124,0] = 9
$data[124,1] = 'Vega Central Mapocho de Santiago'
$data[124,2] = 'Metropolitana'
$data[124,3] = 45140
$data[124,4] = 13
$data[124,5] = 100112005
$data[124,6] = 'Puerro'
$data[124,7] = 'Sin especificar'
$data[124,8] = 'Primera'
$data[124,9] = 70
$data[124,10] = 8000
$data[124,11] = 8000
$data[124,12] = 8000
$data[124,13] = '$/paquete 20 unidades'
$data[124,14] = 'Provincia de Chacabuco'
$data[124,15] = 400
$data[124,16] = 20
$data[124,17] = 'Hortaliza'
$data[125,0] = 9
$data[125,1] = 'Vega Central Mapocho de Santiago'
$data[125,2] = 'Metropolitana'
$data[125,3] = 45119
$data[125,4] = 13
$data[125,5] = 100112005
$data[125,6] = 'Puerro'
$data[125,7] = 'Sin especificar'
$data[125,8] = 'Primera'
$data[125,9] = 70
$data[125,10] = 8000
$data[125,11] = 8000
$data[125,12] = 8000
$data[125,13] = '$/paquete 20 unidades'
$data[125,14] = 'Provincia de Chacabuco'
$data[125,15] = 400
$data[125,16] = 20
$data[125,17] = 'Hortaliza'
$data[126,0] = 9
$data[126,1] = 'Vega Central Mapocho de Santiago'
$data[126,2] = 'Metropolitana'
$data[126,3] = 44526
$data[126,4] = 13
$data[126,5] = 100112005
$data[126,6] = 'Puerro'
$data[126,7] = 'Sin especificar'
$data[126,8] = 'Primera'
$data[126,9] = 170
$data[126,10] = 6000
$data[126,11] = 7000
$data[126,12] = 6500
$data[126,13] = '$/paquete 20 unidades'
$data[126,14] = 'Provincia de Chacabuco'
$data[126,15] = 325
$data[126,16] = 20
$data[126,17] = 'Hortaliza'
$data[127,0] = 9
$data[127,1] = 'Vega Central Mapocho de Santiago'
$data[127,2] = 'Metropolitana'
$data[127,3] = 44524
$data[127,4] = 13
$data[127,5] = 100112005
$data[127,6] = 'Puerro'
$data[127,7] = 'Sin especificar'
$data[127,8] = 'Primera'
$data[127,9] = 160
$data[127,10] = 6000
$data[127,11] = 7000
$data[127,12] = 6500
$data[127,13] = '$/paquete 20 unidades'
$data[127,14] = 'Provincia de Chacabuco'
$data[127,15] = 325
$data[127,16] = 20
$data[127,17] = 'Hortaliza'
$data[128,0] = 9
$data[128,1] = 'Vega Central Mapocho de Santiago'
$data[128,2] = 'Metropolitana'
$data[128,3] = 44505
$data[128,4] = 13
$data[128,5] = 100112005
$data[128,6] = 'Puerro'
$data[128,7] = 'Sin especificar'
$data[128,8] = 'Primera'
$data[128,9] = 160
$data[128,10] = 6000
$data[128,11] = 7000
$data[128,12] = 6500
$data[128,13] = '$/paquete 20 unidades'
$data[128,14] = 'Provincia de Chacabuco'
$data[128,15] = 325
$data[128,16] = 20
$data[128,17] = 'Hortaliza'
$data[129,0] = 9
$data[129,1] = 'Vega Central Mapocho de Santiago'
$data[129,2] = 'Metropolitana'
$data[129,3] = 44266
$data[129,4] = 13
$data[129,5] = 100112005
$data[129,6] = 'Puerro'
$data[129,7] = 'Sin especificar'
$data[129,8] = 'Primera'
$data[129,9] = 50
$data[129,10] = 8000
$data[129,11] = 8000
$data[129,12] = 8000
$data[129,13] = '$/paquete 20 unidades'
$data[129,14] = 'Provincia de Chacabuco'
$data[129,15] = 400
$data[129,16] = 20
$data[129,17] = 'Hortaliza'
$data[130,0] = 9
$data[130,1] = 'Vega Central Mapocho de Santiago'
$data[130,2] = 'Metropolitana'
$data[130,3] = 45133
$data[130,4] = 13
$data[130,5] = 100112005
$data[130,6] = 'Puerro'
$data[130,7] = 'Sin especificar'
$data[130,8] = 'Primera'
$data[130,9] = 70
$data[130,10] = 8000
$data[130,11] = 8000
$data[130,12] = 8000
$data[130,13] = '$/paquete 20 unidades'
$data[130,14] = 'Provincia de Chacabuco'
$data[130,15] = 400
$data[130,16] = 20
$data[130,17] = 'Hortaliza'

$ws.Range("A2:R132").Value = $data

# The freshly-appended row 132 needs the same date number format as the rest of column D.
$ws.Cells.Item(132, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

$ws.Range("A1:R132").Select()
